# The commit swaps the two embedded theme parts of this deck:
#   ppt/theme/theme1.xml (the slide master's theme, "Integral")
#   ppt/theme/theme2.xml (the notes master's theme, "Office Theme")
# end up holding each other's content - i.e. the deck's slide master
# switches from the "Integral" palette to the stock "Office" palette.
#
# The <a:fontScheme>/<a:fmtScheme> blocks of the two themes are already
# byte-identical, so the only real content delta is the 12-colour
# <a:clrScheme>. We drive that through the Design/Master colour-scheme
# object, which PowerPoint exposes as a flat list of RGBColor slots
# mapped 1:1 onto dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

function Set-SchemeColor($index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # VBA/COM RGBColor.RGB packs the value as &H00BBGGRR.
    $scheme.Colors($index).RGB = ($b * 65536) + ($g * 256) + $r
}

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
Set-SchemeColor 1  "000000"
Set-SchemeColor 2  "FFFFFF"
Set-SchemeColor 3  "44546A"
Set-SchemeColor 4  "E7E6E6"
Set-SchemeColor 5  "5B9BD5"
Set-SchemeColor 6  "ED7D31"
Set-SchemeColor 7  "A5A5A5"
Set-SchemeColor 8  "FFC000"
Set-SchemeColor 9  "4472C4"
Set-SchemeColor 10 "70AD47"
Set-SchemeColor 11 "0563C1"
Set-SchemeColor 12 "954F72"
